# Generate Report for Archive
# - Localization status moves from "Ready for handoff" to "In Translation"
#   (shared string used on the Overview sheet's zh-cn/de-de status cells and
#   on each locale sheet's own Status column).
# - The Status column is narrower now that "In Translation" is shorter than
#   "Ready for handoff", so the report's column widths are re-tightened to
#   match the new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text everywhere it appears.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Re-tighten the Status columns to fit the shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.6
$wsOverview.Columns.Item(6).ColumnWidth = 12.6
$wsZhCn.Columns.Item(3).ColumnWidth = 12.6
$wsDeDe.Columns.Item(3).ColumnWidth = 12.6
